$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.869.15"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").Value = "'2.032.41"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'227.79"
$ws.Range("E5").Value = "  -0.95%  "

$ws.Range("D6").Value = "'0.612"
$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("D7").Value = "'60.26"
$ws.Range("E7").Value = "  +7.28%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("D10").Value = "'0.0811"
$ws.Range("E10").Value = "  +1.22%  "

$ws.Range("E11").Value = "  +0.91%  "

$ws.Range("E12").Value = "  +1.27%  "

$ws.Range("D13").Value = "'2.333.55"
$ws.Range("E13").Value = "  +0.17%  "

$ws.Range("D14").Value = "'21.16"
$ws.Range("E14").Value = "  +4.28%  "

$ws.Range("D15").Value = "'0.755"
$ws.Range("E15").Value = "  +1.51%  "

$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("D17").Value = "'2.035.35"
$ws.Range("E17").Value = "  +0.33%  "

$ws.Range("D18").Value = "'37.820.24"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("E19").Value = "  -2.13%  "

$ws.Range("D20").Value = "'69.57"
$ws.Range("E20").Value = "  +0.68%  "

$ws.Range("D21").Value = "'0.0₃0825"
$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").Value = "'224.38"
$ws.Range("E22").Value = "  +0.39%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("E24").Value = "  -0.94%  "

$ws.Range("E25").Value = "  -2.96%  "

$ws.Range("D26").Value = "'165.03"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").Value = "'9.16"
$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("E28").Value = "  -2.50%  "

$ws.Range("D29").Value = "'18.89"
$ws.Range("E29").Value = "  +0.41%  "

$ws.Range("D30").Value = "'1.29"
$ws.Range("E30").Value = "  -3.16%  "

$ws.Range("D31").Value = "'0.120"
$ws.Range("E31").Value = "  +1.67%  "

$ws.Range("E32").Value = "  -1.15%  "

$ws.Range("E33").Value = "  +1.97%  "

$ws.Range("E34").Value = "  +0.54%  "

$ws.Range("E35").Value = "  -0.70%  "

$ws.Range("D36").Value = "'6.32"
$ws.Range("E36").Value = "  +9.38%  "

$ws.Range("D37").Value = "'2.26"
$ws.Range("E37").Value = "  -2.87%  "

$ws.Range("D38").Value = "'3.23"
$ws.Range("E38").Value = "  -0.85%  "

$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("D40").Value = "'1.529.71"
$ws.Range("E40").Value = "  +3.57%  "

$ws.Range("D41").Value = "'0.0217"
$ws.Range("E41").Value = "  +0.90%  "

$ws.Range("D42").Value = "'97.15"
$ws.Range("E42").Value = "  +1.91%  "

$ws.Range("D43").Value = "'16.57"
$ws.Range("E43").Value = "  +0.60%  "

$ws.Range("D44").Value = "'0.0918"
$ws.Range("E44").Value = "  -1.73%  "

$ws.Range("E45").Value = "  -1.33%  "

$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("E47").Value = "  -5.92%  "

$ws.Range("E48").Value = "  +0.62%  "

$ws.Range("E49").Value = "  -0.71%  "

$ws.Range("E50").Value = "  -1.07%  "

$ws.Range("D51").Value = "'2.222.47"
$ws.Range("E51").Value = "  +0.12%  "
